# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 07:39"

# Row 51 - Israel
$ws.Range("B51").Value = 19894
$ws.Range("C51").Value = 111
$ws.Range("D51").Value = 15499
$ws.Range("E51").Value = 4092

# Row 76 - Uzbekistan
$ws.Range("B76").Value = 5697
$ws.Range("C76").Value = 15
$ws.Range("E76").Value = 1547

# Row 87 - El Salvador
$ws.Range("D87").Value = 2152
$ws.Range("E87").Value = 1832
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 82

# Row 94 - Tailandia
$ws.Range("B94").Value = 3141
$ws.Range("C94").Value = 6
$ws.Range("D94").Value = 2997
$ws.Range("E94").Value = 86

# Row 160 - Birmania
$ws.Range("B160").Value = 263
$ws.Range("C160").Value = 1
$ws.Range("E160").Value = 72

# Row 185 - Butan
$ws.Range("D185").Value = 25
$ws.Range("E185").Value = 42

# Row 202 - Dominica
$ws.Range("D202").Value = 18
$ws.Range("E202").Value = 0
